$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1157, "amigo", "aaaa", 3, 4),
    @(1158, "anonimo", "sss", 3, 4),
    @(1160, "bonito", "muy bonito", 4, 5),
    @(1161, "masbonito", "muy bonito", 4, 5),
    @(1162, "aunmasbonito", "muy bonito", 4, 5),
    @(1163, "superbonito", "muy bonito", 4, 5),
    @(1164, "superbonito2", "muy bonito", 4, 5),
    @(1165, "pepino", "", 0, 0)
)

$startRow = 203
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
}
